$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.180.87"
$ws.Range("E2").Value = "  +2.81%  "
$ws.Range("D3").Value = "1.810.74"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.31"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3929"
$ws.Range("E7").Value = "  +3.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3484"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.45"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.188"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07547"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9998"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.07"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.510"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "1.813.59"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.145"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001104"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06692"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.91"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.76"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.577"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").Value = "28.179.79"
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.46"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.404"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.491"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.532"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.32"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.63"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "2.018.57"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.55"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.159"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.015"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08837"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.06"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6942"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06559"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02424"
$ws.Range("E38").Value = "  +3.08%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.457"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.607"
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2210"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.475"
$ws.Range("E43").Value = "  -4.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.55"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6426"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.872"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.143"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.40"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07190"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.97"
$ws.Range("E51").Value = "  +0.35%  "
